$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.851650714874268
$ws.Range("B1").Value = 7.223989009857178
$ws.Range("C1").Value = 5.801651954650879
$ws.Range("D1").Value = 2.229875564575195
$ws.Range("E1").Value = 1.406411647796631
